# Scheduled-runner update: refresh market-board price columns (H-N)
# for the affected leve rows across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 63: Summoning for Dummies / Archaeoskin Codex
$ws.Range("H63").Value = 15000
$ws.Range("I63").Value = 10000
$ws.Range("K63").Value = 10000
$ws.Range("M63").Value = -9376

# Row 66: Summoning the Courage to Be Different (L) / Archaeoskin Codex
$ws.Range("H66").Value = 15000
$ws.Range("I66").Value = 10000
$ws.Range("K66").Value = 30000
$ws.Range("M66").Value = -26880

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 2740.5293
$ws.Range("I135").Value = 2604.9092
$ws.Range("J135").Value = 2989.1667
$ws.Range("K135").Value = 23444.1828
$ws.Range("L135").Value = 26902.5003
$ws.Range("M135").Value = -20909.1828
$ws.Range("N135").Value = -31972.5003

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 28572868
$ws.Range("I137").Value = 41667590
$ws.Range("J137").Value = 2563.818
$ws.Range("K137").Value = 125002770
$ws.Range("L137").Value = 7691.454000000001
$ws.Range("M137").Value = -125000220
$ws.Range("N137").Value = -12791.454

$ws = $wb.Worksheets.Item("ARM")
# Row 35: Need for Mead / Conical Alembic
$ws.Range("H35").Value = 995
$ws.Range("I35").Value = 995
$ws.Range("K35").Value = 995
$ws.Range("M35").Value = -589

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3629.25
$ws.Range("I61").Value = 2574.1333
$ws.Range("J61").Value = 4846.6924
$ws.Range("K61").Value = 2574.1333
$ws.Range("L61").Value = 4846.6924
$ws.Range("M61").Value = -2362.1333
$ws.Range("N61").Value = -5270.6924

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 5610.9287
$ws.Range("I74").Value = 1384.6666
$ws.Range("J74").Value = 18289.715
$ws.Range("K74").Value = 1384.6666
$ws.Range("L74").Value = 18289.715
$ws.Range("M74").Value = -510.6666
$ws.Range("N74").Value = -20037.715

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 5610.9287
$ws.Range("I77").Value = 1384.6666
$ws.Range("J77").Value = 18289.715
$ws.Range("K77").Value = 6923.333000000001
$ws.Range("L77").Value = 91448.575
$ws.Range("M77").Value = -2555.333000000001
$ws.Range("N77").Value = -100184.575

# Row 123: The Armoire Is Open / High Durium Armguards of Maiming
$ws.Range("H123").Value = 31079
$ws.Range("J123").Value = 31079
$ws.Range("L123").Value = 31079
$ws.Range("N123").Value = -40879

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2115
$ws.Range("I132").Value = 1598.4048
$ws.Range("J132").Value = 3058.348
$ws.Range("K132").Value = 4795.2144
$ws.Range("L132").Value = 9175.044
$ws.Range("M132").Value = -2265.2144
$ws.Range("N132").Value = -14235.044

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3629.25
$ws.Range("I136").Value = 2574.1333
$ws.Range("J136").Value = 4846.6924
$ws.Range("K136").Value = 7722.3999
$ws.Range("L136").Value = 14540.0772
$ws.Range("M136").Value = -5172.3999
$ws.Range("N136").Value = -19640.0772

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 28574850
$ws.Range("I134").Value = 50001996
$ws.Range("J134").Value = 5321.2666
$ws.Range("K134").Value = 150005988
$ws.Range("L134").Value = 15963.7998
$ws.Range("M134").Value = -150003453
$ws.Range("N134").Value = -21033.7998

$ws = $wb.Worksheets.Item("CRP")
# Row 5: Bowing Out / Maple Shortbow
$ws.Range("H5").Value = 1803.091
$ws.Range("I5").Value = 75.2
$ws.Range("J5").Value = 3243
$ws.Range("K5").Value = 75.2
$ws.Range("L5").Value = 3243
$ws.Range("M5").Value = 36.8
$ws.Range("N5").Value = -3467

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1762.6875
$ws.Range("I31").Value = 1671.6428
$ws.Range("J31").Value = 2400
$ws.Range("K31").Value = 1671.6428
$ws.Range("L31").Value = 2400
$ws.Range("M31").Value = -1376.6428
$ws.Range("N31").Value = -2990

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1762.6875
$ws.Range("I34").Value = 1671.6428
$ws.Range("J34").Value = 2400
$ws.Range("K34").Value = 1671.6428
$ws.Range("L34").Value = 2400
$ws.Range("M34").Value = -1469.6428
$ws.Range("N34").Value = -2804

# Row 39: An Expected Tourney / Ash Cavalry Bow
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("K39").Value = 1000
$ws.Range("M39").Value = -609

# Row 49: Bend It Like Durendaire / Ash Cavalry Bow
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("K49").Value = 1000
$ws.Range("M49").Value = -818

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2300.303
$ws.Range("I58").Value = 1595.6111
$ws.Range("J58").Value = 3145.9333
$ws.Range("K58").Value = 1595.6111
$ws.Range("L58").Value = 3145.9333
$ws.Range("M58").Value = -1392.6111
$ws.Range("N58").Value = -3551.9333

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 15599.667
$ws.Range("I62").Value = 19649.584
$ws.Range("K62").Value = 19649.584
$ws.Range("M62").Value = -19025.584

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 15599.667
$ws.Range("I65").Value = 19649.584
$ws.Range("K65").Value = 98247.92
$ws.Range("M65").Value = -95127.92

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2391.9143
$ws.Range("I132").Value = 1992.9524
$ws.Range("J132").Value = 2990.3572
$ws.Range("K132").Value = 5978.857199999999
$ws.Range("L132").Value = 8971.071599999999
$ws.Range("M132").Value = -3448.857199999999
$ws.Range("N132").Value = -14031.0716

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 2207.6858
$ws.Range("I134").Value = 1061.8
$ws.Range("K134").Value = 3185.4
$ws.Range("M134").Value = -650.3999999999996

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2300.303
$ws.Range("I136").Value = 1595.6111
$ws.Range("J136").Value = 3145.9333
$ws.Range("K136").Value = 4786.8333
$ws.Range("L136").Value = 9437.7999
$ws.Range("M136").Value = -2236.8333
$ws.Range("N136").Value = -14537.7999

$ws = $wb.Worksheets.Item("CUL")
# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 6315554
$ws.Range("I137").Value = 10002010
$ws.Range("J137").Value = 171460.5
$ws.Range("K137").Value = 30006030
$ws.Range("L137").Value = 514381.5
$ws.Range("M137").Value = -30000930
$ws.Range("N137").Value = -524581.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 7070.5884
$ws.Range("I70").Value = 7264.2856
$ws.Range("J70").Value = 6166.6665
$ws.Range("K70").Value = 7264.2856
$ws.Range("L70").Value = 6166.6665
$ws.Range("M70").Value = -6994.2856
$ws.Range("N70").Value = -6706.6665

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 7070.5884
$ws.Range("I73").Value = 7264.2856
$ws.Range("J73").Value = 6166.6665
$ws.Range("K73").Value = 7264.2856
$ws.Range("L73").Value = 6166.6665
$ws.Range("M73").Value = -6328.2856
$ws.Range("N73").Value = -8038.6665

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2432.4546
$ws.Range("I132").Value = 2227.4688
$ws.Range("J132").Value = 2717.652
$ws.Range("K132").Value = 6682.4064
$ws.Range("L132").Value = 8152.956
$ws.Range("M132").Value = -4152.4064
$ws.Range("N132").Value = -13212.956

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 5783.4443
$ws.Range("I132").Value = 4810.3
$ws.Range("J132").Value = 6999.875
$ws.Range("K132").Value = 14430.9
$ws.Range("L132").Value = 20999.625
$ws.Range("M132").Value = -11900.9
$ws.Range("N132").Value = -26059.625

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 5291.115
$ws.Range("I136").Value = 2840.5386
$ws.Range("J136").Value = 7741.6924
$ws.Range("K136").Value = 8521.6158
$ws.Range("L136").Value = 23225.0772
$ws.Range("M136").Value = -5971.6158
$ws.Range("N136").Value = -28325.0772

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 17860086
$ws.Range("I132").Value = 25002488
$ws.Range("K132").Value = 75007464
$ws.Range("M132").Value = -75004934

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 11146380
$ws.Range("I136").Value = 22290080
$ws.Range("J136").Value = 2679.6667
$ws.Range("K136").Value = 66870240
$ws.Range("L136").Value = 8039.000100000001
$ws.Range("M136").Value = -66867690
$ws.Range("N136").Value = -13139.0001
